$d = $word.ActiveDocument
$rightQuote = [char]0x2019

# ===========================================================================
# Change 1: insert a new, empty, unnumbered paragraph right before the
# "Top 5 Best Sellers by Total Pizza's Sold and 7. Bottom" heading, and make
# that heading's text (and its paragraph mark) bold + size 14pt (28 half-pts).
# ===========================================================================
$headingText = "Top 5 Best Sellers by Total Pizza" + $rightQuote + "s Sold and 7. Bottom "

$findRange = $d.Content
$findRange.Find.Execute($headingText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headingPara = $findRange.Paragraphs.Item(1)
$headingIndex = $headingPara.Index

# The paragraph holding the chart picture sits immediately before the
# heading and already has exactly the plain/unnumbered pPr we want for the
# inserted blank paragraph, so split right after it.
$prevPara = $d.Paragraphs.Item($headingIndex - 1)
$prevPara.Range.InsertParagraphAfter()

# Re-find the heading (its position shifted by the inserted paragraph) and
# grab a *fresh* Paragraph object via its Index (Range.Paragraphs.Item(1)
# reports the right Start/End/Index, but its own .Range.Text is unreliable,
# so re-fetch from the document's Paragraphs collection).
$findRange2 = $d.Content
$findRange2.Find.Execute($headingText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headingIndex2 = $findRange2.Paragraphs.Item(1).Index
$headingPara2 = $d.Paragraphs.Item($headingIndex2)

# Format the whole paragraph (text run + trailing paragraph mark) so both
# the run's rPr and the paragraph mark's rPr pick up the new formatting.
$headingRange = $headingPara2.Range
$headingRange.Font.Bold = 1
$headingRange.Font.BoldBi = 1
$headingRange.Font.Size = 14
$headingRange.Font.SizeBi = 14

# ===========================================================================
# Change 2: move the "_GoBack" bookmark from the end of the document (right
# after the blank paragraph that used to follow the "BOTTOM 5 WORST SELLERS"
# picture) to immediately after the
# "-- Top 5 best sellers by total pizza Sold" SQL comment text.
# ===========================================================================
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$commentText = "-- Top 5 best sellers by total pizza Sold"
$commentRange = $d.Content
$commentRange.Find.Execute($commentText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# Adding a bookmark with Start == End exactly at this paragraph-end offset
# is mishandled by the host (it silently relocates to the wrong spot), so
# work around it: insert a throwaway character right after the comment,
# wrap *that* (non-collapsed) range with the bookmark, then delete just the
# throwaway character's text. The bookmark range collapses back down to the
# correct, empty position and survives the deletion.
$insertPos = $d.Range($commentRange.End, $commentRange.End)
$insertPos.InsertAfter("X")

$markerRange = $d.Range($commentRange.End, $commentRange.End + 1)
$d.Bookmarks.Add("_GoBack", $markerRange)

$markerRange2 = $d.Range($commentRange.End, $commentRange.End + 1)
$markerRange2.Text = ""
